$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits right after
#    "Visualiza el informe de daños" (it should no longer be there).
#    Bookmarks.Delete() is not reliable in this runtime, so we rewrite the
#    paragraph's content (same text/formatting, without the bookmark tags)
#    using Range.InsertXML on the exact run-range (this swaps the whole
#    paragraph but keeps its pPr/rPr/text identical).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Visualiza el informe de daños", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml1 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0015235F" w:rsidRPr="005D5365" w:rsidRDefault="0040764D" w:rsidP="005D5365"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Visualiza el informe de da&#241;os</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Replace the text "Crear lote" with "Ver informe de daños" and add the
#    "_GoBack" bookmark right after the new text (this mirrors where the
#    bookmark moved to in the edited document).
# ---------------------------------------------------------------------------
$d2 = $word.ActiveDocument
$r2 = $d2.Content
$r2.Find.Execute("Crear lote", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml2 = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005D5365" w:rsidRDefault="00C62783" w:rsidP="007A56CA"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>Ver informe de da&#241;os</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)

Write-Host "Edits applied"
